# Generate Report for Handoff
#
# The "86b20423-...md" file moves into the row-2 slot (status: Handed back,
# in sync with en-US) and the "35efe67b-...md" file moves into the row-3
# slot, now reporting a new status "Ready for handoff" with an updated
# handoff datetime. This touches the Overview sheet plus the per-locale
# (zh-cn / de-de) detail sheets, including the hyperlink display text.

$wb = $excel.ActiveWorkbook

function Set-LinkDisplay {
    param($ws, [string]$addr, [string]$text)
    foreach ($h in $ws.Hyperlinks) {
        if ($h.Range.Address() -eq $addr) {
            $h.TextToDisplay = $text
        }
    }
}

# ---------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------
$ov = $wb.Worksheets.Item("Overview")

$ov.Range("A2").Value = "86b20423-8bb9-4abc-9c66-0eed75ac83da.md"
$ov.Range("B2").Value = "Handed back: in sync with en-US"
$ov.Range("C2").Value = "Handed back: in sync with en-US"
$ov.Range("D2").Value = "2016-35-20 18:35:24"

$ov.Range("A3").Value = "35efe67b-7699-461d-923e-8925f6541628.md"
$ov.Range("B3").Value = "Ready for handoff"
$ov.Range("C3").Value = "Ready for handoff"
$ov.Range("D3").Value = "2016-36-20 18:36:27"

Set-LinkDisplay $ov '$A$2' "86b20423-8bb9-4abc-9c66-0eed75ac83da.md"
Set-LinkDisplay $ov '$A$3' "35efe67b-7699-461d-923e-8925f6541628.md"

# ---------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------
$zh = $wb.Worksheets.Item("zh-cn")

$zh.Range("A2").Value = "86b20423-8bb9-4abc-9c66-0eed75ac83da.md"
$zh.Range("B2").Value = ".md"
$zh.Range("C2").Value = "Handed back: in sync with en-US"
$zh.Range("D2").Value = "86b20423-8bb9-4abc-9c66-0eed75ac83da.f8813fda73210ddd7e97c76b77c62525b5b8a2c9.zh-cn.xlf"
$zh.Range("E2").Value = "2016-03-20 18:35:20"
$zh.Range("F2").Value = "86b20423-8bb9-4abc-9c66-0eed75ac83da.md"
$zh.Range("G2").Value = "86b20423-8bb9-4abc-9c66-0eed75ac83da.f8813fda73210ddd7e97c76b77c62525b5b8a2c9.zh-cn.xlf"
$zh.Range("H2").Value = "2016-03-20 18:35:46"
$zh.Range("I2").Value = "Include"

$zh.Range("A3").Value = "35efe67b-7699-461d-923e-8925f6541628.md"
$zh.Range("B3").Value = ".md"
$zh.Range("C3").Value = "Ready for handoff"
$zh.Range("D3").Value = "35efe67b-7699-461d-923e-8925f6541628.c91227dc15ef442992aeeaf68e97009782a86854.zh-cn.xlf"
$zh.Range("E3").Value = "2016-03-20 18:36:24"
$zh.Range("F3").Value = "35efe67b-7699-461d-923e-8925f6541628.md"
$zh.Range("G3").Value = "35efe67b-7699-461d-923e-8925f6541628.c91227dc15ef442992aeeaf68e97009782a86854.zh-cn.xlf"
$zh.Range("H3").Value = "2016-03-20 18:35:46"
$zh.Range("I3").Value = "Include"

Set-LinkDisplay $zh '$A$2' "86b20423-8bb9-4abc-9c66-0eed75ac83da.md"
Set-LinkDisplay $zh '$B$2' ".md"
Set-LinkDisplay $zh '$D$2' "86b20423-8bb9-4abc-9c66-0eed75ac83da.f8813fda73210ddd7e97c76b77c62525b5b8a2c9.zh-cn.xlf"
Set-LinkDisplay $zh '$F$2' "86b20423-8bb9-4abc-9c66-0eed75ac83da.md"
Set-LinkDisplay $zh '$G$2' "86b20423-8bb9-4abc-9c66-0eed75ac83da.f8813fda73210ddd7e97c76b77c62525b5b8a2c9.zh-cn.xlf"

Set-LinkDisplay $zh '$A$3' "35efe67b-7699-461d-923e-8925f6541628.md"
Set-LinkDisplay $zh '$B$3' ".md"
Set-LinkDisplay $zh '$D$3' "35efe67b-7699-461d-923e-8925f6541628.c91227dc15ef442992aeeaf68e97009782a86854.zh-cn.xlf"
Set-LinkDisplay $zh '$F$3' "35efe67b-7699-461d-923e-8925f6541628.md"
Set-LinkDisplay $zh '$G$3' "35efe67b-7699-461d-923e-8925f6541628.c91227dc15ef442992aeeaf68e97009782a86854.zh-cn.xlf"

# ---------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------
$de = $wb.Worksheets.Item("de-de")

$de.Range("A2").Value = "86b20423-8bb9-4abc-9c66-0eed75ac83da.md"
$de.Range("B2").Value = ".md"
$de.Range("C2").Value = "Handed back: in sync with en-US"
$de.Range("D2").Value = "86b20423-8bb9-4abc-9c66-0eed75ac83da.f8813fda73210ddd7e97c76b77c62525b5b8a2c9.de-de.xlf"
$de.Range("E2").Value = "2016-03-20 18:35:24"
$de.Range("F2").Value = "86b20423-8bb9-4abc-9c66-0eed75ac83da.md"
$de.Range("G2").Value = "86b20423-8bb9-4abc-9c66-0eed75ac83da.f8813fda73210ddd7e97c76b77c62525b5b8a2c9.de-de.xlf"
$de.Range("H2").Value = "2016-03-20 18:35:52"
$de.Range("I2").Value = "Include"

$de.Range("A3").Value = "35efe67b-7699-461d-923e-8925f6541628.md"
$de.Range("B3").Value = ".md"
$de.Range("C3").Value = "Ready for handoff"
$de.Range("D3").Value = "35efe67b-7699-461d-923e-8925f6541628.c91227dc15ef442992aeeaf68e97009782a86854.de-de.xlf"
$de.Range("E3").Value = "2016-03-20 18:36:27"
$de.Range("F3").Value = "35efe67b-7699-461d-923e-8925f6541628.md"
$de.Range("G3").Value = "35efe67b-7699-461d-923e-8925f6541628.c91227dc15ef442992aeeaf68e97009782a86854.de-de.xlf"
$de.Range("H3").Value = "2016-03-20 18:35:52"
$de.Range("I3").Value = "Include"

Set-LinkDisplay $de '$A$2' "86b20423-8bb9-4abc-9c66-0eed75ac83da.md"
Set-LinkDisplay $de '$B$2' ".md"
Set-LinkDisplay $de '$D$2' "86b20423-8bb9-4abc-9c66-0eed75ac83da.f8813fda73210ddd7e97c76b77c62525b5b8a2c9.de-de.xlf"
Set-LinkDisplay $de '$F$2' "86b20423-8bb9-4abc-9c66-0eed75ac83da.md"
Set-LinkDisplay $de '$G$2' "86b20423-8bb9-4abc-9c66-0eed75ac83da.f8813fda73210ddd7e97c76b77c62525b5b8a2c9.de-de.xlf"

Set-LinkDisplay $de '$A$3' "35efe67b-7699-461d-923e-8925f6541628.md"
Set-LinkDisplay $de '$B$3' ".md"
Set-LinkDisplay $de '$D$3' "35efe67b-7699-461d-923e-8925f6541628.c91227dc15ef442992aeeaf68e97009782a86854.de-de.xlf"
Set-LinkDisplay $de '$F$3' "35efe67b-7699-461d-923e-8925f6541628.md"
Set-LinkDisplay $de '$G$3' "35efe67b-7699-461d-923e-8925f6541628.c91227dc15ef442992aeeaf68e97009782a86854.de-de.xlf"
